$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.129.56'
$ws.Range("D3").Value = '1.833.13'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9992'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '241.39'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.48%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6579'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.05%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07405'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.61%  '
$ws.Range("E9").Value = '  -1.03%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '22.87'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.46%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07741'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.75%  '
$ws.Range("D12").Value = '1.843.10'
$ws.Range("E12").Value = '  +0.41%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.990'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.59%  '
$ws.Range("E14").Value = '  -1.29%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '83.13'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -3.59%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.109'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.19%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000008577'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +4.08%  '
$ws.Range("D18").Value = '29.139.22'
$ws.Range("E18").Value = '  +0.29%  '
$ws.Range("D19").Value = '2.084.50'
$ws.Range("E19").Value = '  -0.07%  '
$ws.Range("E20").Value = '  -1.06%  '
$ws.Range("E21").Value = '  -0.65%  '
$ws.Range("E22").Value = '  +0.21%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.106'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.64%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.0000'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.00%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '161.27'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.13%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1404'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.02%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.572'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.26%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '17.98'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.30%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.508'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.38%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.099'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -3.14%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.041'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.46%  '
$ws.Range("E32").Value = '  -1.20%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05277'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.63%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.865'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.32%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7373'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.25%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.650'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.11%  '
$ws.Range("D38").Value = '1.299.70'
$ws.Range("E38").Value = '  -1.29%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01785'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.13%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.742'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.96%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9136'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.93%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.034'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.79%  '
$ws.Range("E43").Value = '  +0.08%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.08157'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +9.23%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '102.08'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.29%  '
$ws.Range("D46").Value = '1.986.87'
$ws.Range("E46").Value = '  +0.07%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5131'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.74%  '
$ws.Range("B48").Value = 'Aave'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '63.69'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.08%  '
$ws.Range("B49").Value = 'RenderToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.747'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.76%  '
$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.05845'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.43%  '
$ws.Range("B51").Value = 'Aptos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '6.760'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.03%  '
